$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2:O73").Value = "2022-08-09 20:57:24"
$ws.Range("N56").Value = "Naturaline Damen String schwarz L 9.95 Schweizer Franken"
